# Renumber the TestcaseID column (A2:A5) on Sheet1:
#   A2: TC01 -> TC05
#   A3: TC02 -> TC03
#   A4: TC03 -> TC01
#   A5: TC04 -> TC02
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "TC05"
$ws.Range("A3").Value = "TC03"
$ws.Range("A4").Value = "TC01"
$ws.Range("A5").Value = "TC02"

# Update the active selection to F2 (was J3)
[void]$ws.Range("F2").Select()
